# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.52 = 39227.26 pesos`n✅ 39227.26 pesos = 9.46 = 967.29 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 105.08
$wsTasas.Range("O10").Value = 4122
$wsTasas.Range("N12").Value = 4145
$wsTasas.Range("O12").Value = 102.21
